# Refresh the cryptocurrency price/volume snapshot in the "cryptos" sheet.
# - Columns B/C (coin name/link) shift down one row for rows 8-21 because a new
#   coin, "LidoStakedEther", is inserted at row 8 (its row 22 predecessor,
#   "WrappedliquidstakedEther2.0", drops off the bottom of the tracked list).
# - Columns D/E (price, 1h volume) are refreshed for every data row (2-51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "29.845.10"; E = "  -0.79%  " },
    @{ Row = 3; D = "1.894.09"; E = "  -0.26%  " },
    @{ Row = 4; D = "0.9997"; E = "  -0.12%  " },
    @{ Row = 5; D = "0.7600"; E = "  +4.18%  " },
    @{ Row = 6; D = "239.66"; E = "  -1.25%  " },
    @{ Row = 7; D = "0.9996"; E = "  -0.20%  " },
    @{ Row = 8; B = "LidoStakedEther"; C = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"; D = "1.891.63"; E = "  +0.14%  " },
    @{ Row = 9; B = "Cardano"; C = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"; D = "0.3056"; E = "  -1.33%  " },
    @{ Row = 10; B = "Solana"; C = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"; D = "25.52"; E = "  -2.49%  " },
    @{ Row = 11; B = "Dogecoin"; C = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"; D = "0.06829"; E = "  -0.79%  " },
    @{ Row = 12; B = "TRON"; C = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"; D = "0.07947"; E = "  +0.09%  " },
    @{ Row = 13; B = "Polygon"; C = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"; D = "0.7434"; E = "  -3.56%  " },
    @{ Row = 14; B = "WrappedEther"; C = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; D = "1.892.30"; E = "  -0.10%  " },
    @{ Row = 15; B = "Polkadot"; C = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"; D = "5.160"; E = "  -1.57%  " },
    @{ Row = 16; B = "Litecoin"; C = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"; D = "90.68"; E = "  -0.27%  " },
    @{ Row = 17; B = "WrappedBTC"; C = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"; D = "29.842.69"; E = "  -0.75%  " },
    @{ Row = 18; B = "Avalanche"; C = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"; D = "13.94"; E = "  -1.43%  " },
    @{ Row = 19; B = "Uniswap"; C = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"; D = "5.949"; E = "  +3.04%  " },
    @{ Row = 20; B = "BitcoinCash"; C = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"; D = "241.81"; E = "  +1.92%  " },
    @{ Row = 21; B = "ShibaInu"; C = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"; D = "0.000007676"; E = "  -0.99%  " },
    @{ Row = 22; B = "Dai"; C = "https://coinranking.com/coin/MoTuySvg7+dai-dai"; D = "0.9992"; E = "  -0.22%  " },
    @{ Row = 23; D = "1.000"; E = "  -0.12%  " },
    @{ Row = 24; D = "6.950"; E = "  +0.57%  " },
    @{ Row = 25; D = "167.06"; E = "  +0.89%  " },
    @{ Row = 26; D = "9.219"; E = "  -0.83%  " },
    @{ Row = 27; D = "18.68"; E = "  -1.21%  " },
    @{ Row = 28; D = "0.1292"; E = "  +1.81%  " },
    @{ Row = 29; D = "2.025"; E = "  +0.49%  " },
    @{ Row = 30; D = "1.387"; E = "  +2.24%  " },
    @{ Row = 31; D = "1.513"; E = "  -1.49%  " },
    @{ Row = 32; D = "4.249"; E = "  -0.99%  " },
    @{ Row = 33; D = "4.045"; E = "  -0.60%  " },
    @{ Row = 34; D = "0.05202"; E = "  +2.32%  " },
    @{ Row = 35; D = "1.253"; E = "  -1.37%  " },
    @{ Row = 36; D = "0.7262"; E = "  -1.06%  " },
    @{ Row = 37; D = "2.709"; E = "  -0.82%  " },
    @{ Row = 38; D = "0.01915"; E = "  -0.16%  " },
    @{ Row = 39; D = "2.773"; E = "  +0.16%  " },
    @{ Row = 40; D = "6.144"; E = "  -3.11%  " },
    @{ Row = 41; D = "0.4401"; E = "  -0.54%  " },
    @{ Row = 42; D = "71.66"; E = "  -3.94%  " },
    @{ Row = 43; D = "0.9992" },
    @{ Row = 44; D = "1.885"; E = "  -1.94%  " },
    @{ Row = 45; D = "0.8275"; E = "  -0.97%  " },
    @{ Row = 46; D = "7.608"; E = "  +0.37%  " },
    @{ Row = 47; D = "99.92"; E = "  -0.86%  " },
    @{ Row = 48; D = "9.746"; E = "  +0.22%  " },
    @{ Row = 49; D = "2.044.98"; E = "  -0.19%  " },
    @{ Row = 50; D = "35.96"; E = "  -4.60%  " },
    @{ Row = 51; D = "0.05940"; E = "  -0.27%  " }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("B")) { $ws.Cells.Item($u.Row, 2).Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Cells.Item($u.Row, 3).Value = $u.C }
    if ($u.ContainsKey("D")) {
        # Force text storage so numeric-looking strings (e.g. "1.000",
        # "0.9997") keep their literal digits/trailing zeros instead of
        # being parsed into a Double by the COM value setter.
        $cell = $ws.Cells.Item($u.Row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    if ($u.ContainsKey("E")) {
        $cell = $ws.Cells.Item($u.Row, 5)
        $cell.NumberFormat = "@"
        $cell.Value = $u.E
        $cell.Style = "Normal"
    }
}
